$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 31; existing rows 31-33 shift down to 32-34.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new weekly record.
$ws.Cells.Item(31, 1).Value = 8
$ws.Cells.Item(31, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(31, 3).Value = "Coquimbo"
$ws.Cells.Item(31, 4).Value = 45132
$ws.Cells.Item(31, 5).Value = 4
$ws.Cells.Item(31, 6).Value = 100112013
$ws.Cells.Item(31, 7).Value = "Alcachofa"
$ws.Cells.Item(31, 8).Value = "Española"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 420
$ws.Cells.Item(31, 11).Value = 13000
$ws.Cells.Item(31, 12).Value = 14000
$ws.Cells.Item(31, 13).Value = 13500
$ws.Cells.Item(31, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(31, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value = 450
$ws.Cells.Item(31, 17).Value = 30
$ws.Cells.Item(31, 18).Value = "Hortaliza"
